# Weekly data refresh for the "Apio" sheet.
# A new observation is inserted at row 30 (new date 44558, volume 100),
# which pushes every existing observation down by one row (rows 30-162
# become rows 31-163); the record that used to be the very last one
# (row 162) becomes the new row 163.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstDataRow = 30
$lastDataRow  = 162
$newLastRow   = $lastDataRow + 1

# Columns whose values actually vary from one observation to the next.
# (A,B,C,E,F,G,I,N,Q,R are constant for every row in this sheet.)
$cols = @(4, 8, 10, 11, 12, 13, 15, 16)   # D,H,J,K,L,M,O,P

# 1) Create the brand-new last row (163) as a copy of the current last
#    row (162) -- it is about to be overwritten by the shift below.
foreach ($col in (1..18)) {
    $ws.Cells.Item($newLastRow, $col).Value = $ws.Cells.Item($lastDataRow, $col).Value()
}

# 2) Shift every observation down by one row: row r takes what used to
#    be in row r-1. Walk from the bottom up so we never clobber a value
#    before we have read it.
for ($r = $lastDataRow; $r -ge ($firstDataRow + 1); $r--) {
    foreach ($col in $cols) {
        $srcValue = $ws.Cells.Item($r - 1, $col).Value()
        $ws.Cells.Item($r, $col).Value = $srcValue
    }
}

# 3) Write the brand-new observation into row 30 itself. Only the date
#    and volume are genuinely new; price/origin columns are unchanged.
$ws.Range("D30").Value = 44558
$ws.Range("J30").Value = 100
